$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '43.941.15'
$ws.Range('E2').Value = '  +0.48%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.365.16'
$ws.Range('E3').Value = '  +2.96%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range('D5') '0.661'
$ws.Range('E5').Value = '  +2.73%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D6') '235.46'
$ws.Range('E6').Value = '  +1.58%  '

# Row 7
Set-TextValue $ws.Range('D7') '73.00'
$ws.Range('E7').Value = '  +13.28%  '

# Row 8
$ws.Range('E8').Value = '  +0.13%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.539'
$ws.Range('E9').Value = '  +22.24%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.0989'
$ws.Range('E10').Value = '  +2.41%  '

# Row 11
Set-TextValue $ws.Range('D11') '28.39'
$ws.Range('E11').Value = '  +6.69%  '

# Row 12
Set-TextValue $ws.Range('D12') '2.732.95'
$ws.Range('E12').Value = '  +4.38%  '

# Row 13
$ws.Range('E13').Value = '  +2.28%  '

# Row 14
Set-TextValue $ws.Range('D14') '16.96'
$ws.Range('E14').Value = '  +11.71%  '

# Row 15
Set-TextValue $ws.Range('D15') '6.66'
$ws.Range('E15').Value = '  +9.22%  '

# Row 16
Set-TextValue $ws.Range('D16') '0.885'
$ws.Range('E16').Value = '  +6.56%  '

# Row 17
Set-TextValue $ws.Range('D17') '2.368.22'
$ws.Range('E17').Value = '  +3.54%  '

# Row 18
Set-TextValue $ws.Range('D18') '43.858.25'
$ws.Range('E18').Value = '  +0.58%  '

# Row 19
$ws.Range('E19').Value = '  +3.83%  '

# Row 20
Set-TextValue $ws.Range('D20') '76.37'
$ws.Range('E20').Value = '  +4.21%  '

# Row 21
$ws.Range('E21').Value = '  +3.72%  '

# Row 22
Set-TextValue $ws.Range('D22') '251.81'
$ws.Range('E22').Value = '  +1.41%  '

# Row 23
$ws.Range('E23').Value = '  +2.28%  '

# Row 24
$ws.Range('E24').Value = '  -0.14%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.49'
$ws.Range('E25').Value = '  +2.13%  '

# Row 26
$ws.Range('E26').Value = '  +6.25%  '

# Row 27
$ws.Range('E27').Value = '  -1.41%  '

# Row 28
Set-TextValue $ws.Range('D28') '22.55'
$ws.Range('E28').Value = '  +2.53%  '

# Row 29
Set-TextValue $ws.Range('D29') '173.19'
$ws.Range('E29').Value = '  -0.52%  '

# Row 30
Set-TextValue $ws.Range('D30') '1.56'
$ws.Range('E30').Value = '  +9.42%  '

# Row 31
$ws.Range('E31').Value = '  +1.56%  '

# Row 32
$ws.Range('E32').Value = '  +4.90%  '

# Row 33
$ws.Range('E33').Value = '  +5.32%  '

# Row 34
$ws.Range('E34').Value = '  +4.04%  '

# Row 35
Set-TextValue $ws.Range('D35') '5.13'
$ws.Range('E35').Value = '  +3.62%  '

# Row 36
$ws.Range('E36').Value = '  +4.69%  '

# Row 37
$ws.Range('E37').Value = '  +5.28%  '

# Row 38
$ws.Range('E38').Value = '  -0.62%  '

# Row 39
$ws.Range('E39').Value = '  +6.08%  '

# Row 40
Set-TextValue $ws.Range('D40') '19.58'
$ws.Range('E40').Value = '  +13.87%  '

# Row 41
$ws.Range('E41').Value = '  +1.88%  '

# Row 42
$ws.Range('E42').Value = '  -0.16%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.18'
$ws.Range('E43').Value = '  +8.14%  '

# Row 44
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D44') '0.0971'
$ws.Range('E44').Value = '  +2.83%  '

# Row 45
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D45') '1.22'
$ws.Range('E45').Value = '  +2.54%  '

# Row 46
Set-TextValue $ws.Range('D46') '98.47'
$ws.Range('E46').Value = '  +1.14%  '

# Row 47
Set-TextValue $ws.Range('D47') '4.44'
$ws.Range('E47').Value = '  +0.11%  '

# Row 48
$ws.Range('E48').Value = '  +13.48%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D49') '2.33'
$ws.Range('E49').Value = '  +2.57%  '

# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D50') '1.442.35'
$ws.Range('E50').Value = '  +0.66%  '

# Row 51
Set-TextValue $ws.Range('D51') '2.591.18'
$ws.Range('E51').Value = '  +3.67%  '
